# Update row 8 with the corrected/merged record, then remove the old row 9
# (its data has been folded into the new row 8), matching the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "James "
$ws.Range("B8").Value = "Joy"
$ws.Range("C8").Value = "Thiruvalla"
$ws.Range("D8").Value = "adam@yourmail.com"
$ws.Range("E8").Value = 7894561230
$ws.Range("F8").Value = 33101
$ws.Range("G8").Value = "Jino"
$ws.Range("H8").Value = "Duplicate email entry found,Role is not valid"

# Row 9 is no longer needed; delete it and shift things up.
$ws.Range("A9:H9").EntireRow.Delete()

# Column H needs to widen to fit the longer "Result" text (best-fit recalculation).
$ws.Columns.Item(8).ColumnWidth = 40.61328125
